# Weekly update: a new price record (week of 2021-09-29, serial 44468) is
# inserted as row 50 on the "Locoto" sheet. All subsequent rows (old 50-70)
# shift down by one (new 51-71), and the sheet dimension grows to A1:R71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 50, pushing existing rows 50-70 down to 51-71.
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new price record.
$ws.Cells.Item(50, 1).Value  = 1
$ws.Cells.Item(50, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value  = 44468
$ws.Cells.Item(50, 5).Value  = 15
$ws.Cells.Item(50, 6).Value  = 100112042
$ws.Cells.Item(50, 7).Value  = "Locoto"
$ws.Cells.Item(50, 8).Value  = "Sin especificar"
$ws.Cells.Item(50, 9).Value  = "Primera"
$ws.Cells.Item(50, 10).Value = 120
$ws.Cells.Item(50, 11).Value = 27000
$ws.Cells.Item(50, 12).Value = 28000
$ws.Cells.Item(50, 13).Value = 27500
$ws.Cells.Item(50, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(50, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value = 1375
$ws.Cells.Item(50, 17).Value = 20
$ws.Cells.Item(50, 18).Value = "Hortaliza"
